$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new day's log entry (row 4) ---
# Copy the date-cell formatting from row 2 (A2) onto A4, and the plain
# "session" formatting from B2:D2 onto B4:D4 before writing the values,
# so the new row matches the look of the existing entries.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B2:D2").Copy()
$ws.Range("B4:D4").PasteSpecial(-4122)

$ws.Range("A4").Value = 45223
$ws.Range("B4").Value = 67
$ws.Range("C4").Value = "Created authRepo, authController, firebase instance providers, utilities having imagePIcker and showSnackBar, created and worked on SignUpScreen and setup functions."
$ws.Range("D4").Value = 1

# --- Update the sheet's current selection to cover the logged rows ---
$ws.Range("A1:D4").Select()
